$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 112.71429
$ws.Range("I5").Value = 112.71429
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 112.71429
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 2.285709999999995
$ws.Range("N5").ClearContents()
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H40").Value = 4438.3706
$ws.Range("I40").Value = 4057.4375
$ws.Range("J40").Value = 4992.4546
$ws.Range("K40").Value = 4057.4375
$ws.Range("L40").Value = 4992.4546
$ws.Range("M40").Value = -3882.4375
$ws.Range("N40").Value = -5342.4546
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H113").Value = 5876
$ws.Range("I113").Value = 7490.5713
$ws.Range("K113").Value = 7490.5713
$ws.Range("M113").Value = -4236.5713
$ws.Range("H116").Value = 8998.799999999999
$ws.Range("I116").Value = 10148.5
$ws.Range("J116").Value = 4400
$ws.Range("K116").Value = 10148.5
$ws.Range("L116").Value = 4400
$ws.Range("M116").Value = -6706.5
$ws.Range("N116").Value = -11284
$ws.Range("H138").Value = 2370.8635
$ws.Range("J138").Value = 2894.3333
$ws.Range("L138").Value = 8682.999899999999
$ws.Range("N138").Value = -18962.9999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100.2
$ws.Range("I4").Value = 100.333336
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 100.333336
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 15.666664
$ws.Range("N4").Value = -332
$ws.Range("H43").Value = 7540085.5
$ws.Range("I43").Value = 10020114
$ws.Range("K43").Value = 10020114
$ws.Range("M43").Value = -10019801
$ws.Range("H61").Value = 4112.6665
$ws.Range("I61").Value = 3714.2856
$ws.Range("K61").Value = 3714.2856
$ws.Range("M61").Value = -3502.2856
$ws.Range("H122").Value = 3227.9167
$ws.Range("I122").Value = 2972.889
$ws.Range("K122").Value = 8918.667000000001
$ws.Range("M122").Value = -6468.667000000001
$ws.Range("H136").Value = 4112.6665
$ws.Range("I136").Value = 3714.2856
$ws.Range("K136").Value = 11142.8568
$ws.Range("M136").Value = -8592.856800000001
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3635.5386
$ws.Range("I20").Value = 4206.6
$ws.Range("J20").Value = 1732
$ws.Range("K20").Value = 4206.6
$ws.Range("L20").Value = 1732
$ws.Range("M20").Value = -3959.6
$ws.Range("N20").Value = -2226
$ws.Range("H64").Value = 799.6667
$ws.Range("I64").Value = 787
$ws.Range("K64").Value = 787
$ws.Range("M64").Value = -562
$ws.Range("H67").Value = 799.6667
$ws.Range("I67").Value = 787
$ws.Range("K67").Value = 787
$ws.Range("M67").Value = -7
$ws.Range("H80").Value = 168.05
$ws.Range("I80").Value = 120
$ws.Range("K80").Value = 120
$ws.Range("M80").Value = 878
$ws.Range("H83").Value = 168.05
$ws.Range("I83").Value = 120
$ws.Range("K83").Value = 600
$ws.Range("M83").Value = 4392
$ws.Range("H86").Value = 4346.5386
$ws.Range("I86").Value = 1751
$ws.Range("J86").Value = 5968.75
$ws.Range("K86").Value = 1751
$ws.Range("L86").Value = 5968.75
$ws.Range("M86").Value = -628
$ws.Range("N86").Value = -8214.75
$ws.Range("H89").Value = 4346.5386
$ws.Range("I89").Value = 1751
$ws.Range("J89").Value = 5968.75
$ws.Range("K89").Value = 8755
$ws.Range("L89").Value = 29843.75
$ws.Range("M89").Value = -3139
$ws.Range("N89").Value = -41075.75
$ws.Range("H102").Value = 10556
$ws.Range("I102").Value = 10556
$ws.Range("K102").Value = 10556
$ws.Range("M102").Value = -7311

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1126.3334
$ws.Range("I3").Value = 750.6667
$ws.Range("J3").Value = 1502
$ws.Range("K3").Value = 750.6667
$ws.Range("L3").Value = 1502
$ws.Range("M3").Value = -637.6667
$ws.Range("N3").Value = -1728
$ws.Range("H7").Value = 136.41667
$ws.Range("I7").Value = 121.28571
$ws.Range("K7").Value = 121.28571
$ws.Range("M7").Value = -8.285709999999995
$ws.Range("H22").Value = 1720
$ws.Range("I22").Value = 891
$ws.Range("J22").Value = 3999.75
$ws.Range("K22").Value = 891
$ws.Range("L22").Value = 3999.75
$ws.Range("M22").Value = -541
$ws.Range("N22").Value = -4699.75
$ws.Range("H93").Value = 15000
$ws.Range("I93").Value = 15000
$ws.Range("K93").Value = 15000
$ws.Range("M93").Value = -13128
$ws.Range("H141").Value = 84084
$ws.Range("J141").Value = 84084
$ws.Range("L141").Value = 84084
$ws.Range("N141").Value = -94444

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 3150.5
$ws.Range("I31").Value = 1301
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 3903
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = -3615
$ws.Range("N31").Value = -15576
$ws.Range("H99").Value = 26500
$ws.Range("I99").Value = 50000
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 150000
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -147754
$ws.Range("N99").Value = -13492
$ws.Range("H107").Value = 900
$ws.Range("J107").Value = 900
$ws.Range("L107").Value = 2700
$ws.Range("N107").Value = -6540
$ws.Range("H131").Value = 1893.909
$ws.Range("I131").Value = 1916.75
$ws.Range("J131").Value = 1833
$ws.Range("K131").Value = 5750.25
$ws.Range("L131").Value = 5499
$ws.Range("M131").Value = -710.25
$ws.Range("N131").Value = -15579
$ws.Range("H132").Value = 1566.4445
$ws.Range("I132").Value = 799
$ws.Range("J132").Value = 1950.1666
$ws.Range("K132").Value = 7191
$ws.Range("L132").Value = 17551.4994
$ws.Range("M132").Value = -4661
$ws.Range("N132").Value = -22611.4994
$ws.Range("H139").Value = 2101.5454
$ws.Range("I139").Value = 1779.6666
$ws.Range("K139").Value = 5338.9998
$ws.Range("M139").Value = -198.9997999999996

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 21000
$ws.Range("J49").Value = 21000
$ws.Range("L49").Value = 21000
$ws.Range("N49").Value = -21368
$ws.Range("H70").Value = 10335.333
$ws.Range("I70").Value = 4669
$ws.Range("K70").Value = 4669
$ws.Range("M70").Value = -4399
$ws.Range("H73").Value = 10335.333
$ws.Range("I73").Value = 4669
$ws.Range("K73").Value = 4669
$ws.Range("M73").Value = -3733
$ws.Range("H122").Value = 4917.7144
$ws.Range("I122").Value = 3361
$ws.Range("K122").Value = 10083
$ws.Range("M122").Value = -7633
$ws.Range("H132").Value = 5098.3335
$ws.Range("I132").Value = 4640.5
$ws.Range("K132").Value = 13921.5
$ws.Range("M132").Value = -11391.5
$ws.Range("H135").Value = 245000
$ws.Range("J135").Value = 245000
$ws.Range("L135").Value = 245000
$ws.Range("N135").Value = -255140

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 917.7273
$ws.Range("I22").Value = 886.875
$ws.Range("K22").Value = 886.875
$ws.Range("M22").Value = -591.875
$ws.Range("H27").Value = 917.7273
$ws.Range("I27").Value = 886.875
$ws.Range("K27").Value = 886.875
$ws.Range("M27").Value = -779.875
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H93").Value = 1356.8
$ws.Range("I93").Value = 1356.8
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1356.8
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -108.8
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 7517.4546
$ws.Range("J100").Value = 8653.333000000001
$ws.Range("L100").Value = 8653.333000000001
$ws.Range("N100").Value = -9735.333000000001
$ws.Range("H132").Value = 3060.5454
$ws.Range("I132").Value = 3018.5557
$ws.Range("J132").Value = 3249.5
$ws.Range("K132").Value = 9055.667099999999
$ws.Range("L132").Value = 9748.5
$ws.Range("M132").Value = -6525.667099999999
$ws.Range("N132").Value = -14808.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 27592.25
$ws.Range("J104").Value = 27592.25
$ws.Range("L104").Value = 27592.25
$ws.Range("N104").Value = -34580.25
